$wb = $excel.ActiveWorkbook

# --- 1. Text change: "Ready for handoff" -> "In Translation" ---
# This shared string shows up on the Overview sheet (E2, F2) and on each
# per-locale sheet's Status column (C2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Column width changes ---
# Overview columns E and F, and column C on each locale sheet, get narrower
# (from ~17.22 chars down to ~13.41 chars). The host quantizes ColumnWidth
# to 1/6-character steps, so 12.5 is the closest achievable setting (lands
# the stored OOXML width right on its target bucket).
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
